# issue #5: property land done
#
# Land sheet ("土地" / sheet1):
#  - Normalize the land-parcel name and register date text (drop the
#    "★" marker, the stray spaces and dashes).
#  - Switch the header row from the Chinese labels to the scrape
#    pipeline's machine-readable column names.
#  - Append the extra pipeline metadata columns I:O (property_category,
#    category, date, legislator_name, legislator_id, source_file,
#    index) on both the header and data row.
#
# Building sheet ("建物" / sheet2):
#  - Same normalization of the stray "★" marker / spacing / dashes in
#    the building name and the register date; headers stay untouched.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Sheet 1: 土地 (land)
# ---------------------------------------------------------------------

# Give the new trailing header cells (I1:O1) the same bold/bordered
# look as the existing header row before filling them in.
$ws1.Range("H1").Copy() | Out-Null
$ws1.Range("I1:O1").PasteSpecial(-4122) | Out-Null

# Row 1 header labels.
$ws1.Range("B1").Value = "name"
$ws1.Range("C1").Value = "area"
$ws1.Range("D1").Value = "share_portion"
$ws1.Range("E1").Value = "owner"
$ws1.Range("F1").Value = "register_date"
$ws1.Range("G1").Value = "register_reason"
$ws1.Range("H1").Value = "acquire_value"
$ws1.Range("I1").Value = "property_category"
$ws1.Range("J1").Value = "category"
$ws1.Range("K1").Value = "date"
$ws1.Range("L1").Value = "legislator_name"
$ws1.Range("M1").Value = "legislator_id"
$ws1.Range("N1").Value = "source_file"
$ws1.Range("O1").Value = "index"

# Row 2 data cleanup.
$ws1.Range("B2").Value = "台南市安南區海東段00450069地號"
$ws1.Range("F2").Value = "87年07月03日"

# Row 2 new trailing metadata columns.
$ws1.Range("I2").Value = "land"
$ws1.Range("J2").Value = "normal"
# "2011-06-02" looks like a date to Excel's input parser, so force it
# in as literal text (the leading apostrophe is just an entry marker,
# it is not stored in the cell) and then drop back to the sheet's
# plain/general look so it matches its neighbouring cells.
$ws1.Range("K2").Value = "'2011-06-02"
$ws1.Range("H2").Copy() | Out-Null
$ws1.Range("K2").PasteSpecial(-4122) | Out-Null
$ws1.Range("L2").Value = "許添財"
$ws1.Range("M2").Value = 639
$ws1.Range("N2").Value = "tmp15921"
$ws1.Range("O2").Value = 14

# ---------------------------------------------------------------------
# Sheet 2: 建物 (building)
# ---------------------------------------------------------------------
$ws2.Range("B2").Value = "台南市安南區海東段01774000建號((信託)）"
$ws2.Range("F2").Value = "87年07月03日"
